# Auto-generated: applies the crypto price/volume table update described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.211.75'
$ws.Range('E2').Value = '  +1.70%  '
$ws.Range('D3').Value = '2.679.00'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'522.20"
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').Value = "'146.43"
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').Value = "'0.575"
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('D9').Value = '2.697.29'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').Value = "'6.46"
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('D12').Value = "'0.339"
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').Value = '3.154.59'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').Value = '60.199.76'
$ws.Range('E15').Value = '  +1.69%  '
$ws.Range('D16').Value = "'21.34"
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').Value = "'0.0000138"
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '2.690.45'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').Value = "'351.07"
$ws.Range('E19').Value = '  -2.17%  '
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('D21').Value = "'10.54"
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = "'6.33"
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  +2.75%  '
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('E26').Value = '  +4.00%  '
$ws.Range('E27').Value = '  +0.36%  '
$ws.Range('D28').Value = "'7.35"
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').Value = '0.0₃0816'
$ws.Range('E29').Value = '  -1.34%  '
$ws.Range('D30').Value = "'6.79"
$ws.Range('E30').Value = '  +4.98%  '
$ws.Range('E31').Value = '  +0.29%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'19.14"
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = "'1.59"
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('D34').Value = "'147.16"
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('D35').Value = "'4.31"
$ws.Range('E35').Value = '  +5.00%  '
$ws.Range('E36').Value = '  +8.07%  '
$ws.Range('D37').Value = "'0.950"
$ws.Range('E37').Value = '  -6.92%  '
$ws.Range('D38').Value = "'0.876"
$ws.Range('E38').Value = '  +1.23%  '
$ws.Range('D39').Value = "'1.52"
$ws.Range('E39').Value = '  +6.70%  '
$ws.Range('D40').Value = "'36.89"
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('D41').Value = "'3.71"
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('D42').Value = "'284.22"
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value = "'0.0990"
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = "'20.00"
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = "'0.996"
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = "'0.611"
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('D47').Value = '2.128.51'
$ws.Range('E47').Value = '  +5.35%  '
$ws.Range('D48').Value = "'4.95"
$ws.Range('E48').Value = '  +2.94%  '
$ws.Range('D49').Value = "'0.0540"
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('D50').Value = "'0.0235"
$ws.Range('E50').Value = '  +1.08%  '
$ws.Range('D51').Value = "'19.36"
$ws.Range('E51').Value = '  +3.61%  '
